$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 14; this shifts the existing rows 14-24
# (A14:M24) down to 15-25, matching the "after" row numbering in the diff.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new leetcoder's data.
$ws.Range("A14").Value = 47373
$ws.Range("B14").Value = "https://leetcode.com/u/th2do/"
$ws.Range("C14").Value = 55
$ws.Range("D14").Value = 607
$ws.Range("F14").Value = 49
$ws.Range("G14").Value = 275

# The row that used to be row 14 (now row 15, rank 52374) also gained a
# new value in column I.
$ws.Range("I15").Value = 156

# Update the selection to match the saved workbook state (F15).
$ws.Range("F15").Select() | Out-Null
